# Applies the backend/backend/data.xlsx update:
#  - Rename header E1 and add new headers F1:M1
#  - Add two more raw rows (6,7) in columns A:D
#  - Add three new rows (8,9,10) of payment-style data in columns E:M
#
# Cells whose text looks like a pure number or a date (e.g. "565",
# "777888 ", "09-10-2024") are forced to Text format first so Excel
# keeps them as strings (matching numberStoredAsText / t="str" cells
# in the source data) instead of silently coercing them into numbers
# or date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ----------------------------------------------------------
$ws.Range("E1").Value = "Date and Time"
$ws.Range("F1").Value = "Amount"
$ws.Range("G1").Value = "Project Id"
$ws.Range("H1").Value = "Account number"
$ws.Range("I1").Value = "PO number"
$ws.Range("J1").Value = "Vendor Name"
$ws.Range("K1").Value = "IFSC code"
$ws.Range("L1").Value = "status"
$ws.Range("M1").Value = "Remarks"

# --- Row 6 (columns A:D) --------------------------------------------------
$ws.Range("A6").Value = "test009"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "854 "
$ws.Range("C6").Value = "gowtham "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "565"

# --- Row 7 (columns A:D) --------------------------------------------------
$ws.Range("A7").Value = "test009 "
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "854 "
$ws.Range("C7").Value = "tester "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "777888 "

# --- Row 8 (columns E:M) --------------------------------------------------
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "09-10-2024"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "00 "
$ws.Range("G8").Value = "P23"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "777888 "
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "33 "
$ws.Range("J8").Value = "tester "
$ws.Range("K8").Value = "test007"
$ws.Range("L8").Value = "Submitted"
$ws.Range("M8").Value = "Testing Reason"

# --- Row 9 (columns E:M) --------------------------------------------------
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "09-10-2024"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "100000"
$ws.Range("G9").Value = "P23"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "777888 "
$ws.Range("I9").Value = "g"
$ws.Range("J9").Value = "tester "
$ws.Range("K9").Value = "test007"
$ws.Range("L9").Value = "Submitted"
$ws.Range("M9").Value = "Testing Reason"

# --- Row 10 (columns E:M) -------------------------------------------------
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "09-10-2024"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "00"
$ws.Range("G10").Value = "P23"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "696"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "854 "
$ws.Range("J10").Value = "tester  "
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = "2003"
$ws.Range("L10").Value = "Submitted"
$ws.Range("M10").Value = "Testing Reason"
